# Update the "POP3 Messages" rectangle on slide 1 so that it reads
# "Message Source" (main run) plus a smaller sub-label
# "(e.g. POP3, Manual)" -- per commit "Pointing out Message Source".

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$shp = $s.Shapes.Item("Rectangle 8")

$tr = $shp.TextFrame.TextRange

# Replace the whole run's text first, keeping the existing run
# formatting (sz=1400, tx2 solid fill) for the "Message Source " part.
$tr.Text = "Message Source (e.g. POP3, Manual)"

# Grab just the appended "(e.g. POP3, Manual)" portion and shrink it to
# 10pt so it reads as a secondary/explanatory label. Setting the font
# size on this sub-range splits the paragraph into two <a:r> runs while
# keeping the rest of the run formatting (language, color) intact.
$prefixLen = "Message Source ".Length
$suffixLen = $tr.Length - $prefixLen
$suffix = $tr.Characters($prefixLen + 1, $suffixLen)
$suffix.Font.Size = 10
